$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix floating point number formatting in column H (Importe) ---
# The amounts are stored as text (e.g. "4.850,00"); the fix converts them
# to "4850.00" (drop thousands "." separator, "," decimal -> ".").
$importeRange = $ws.Range("H2:H231")
$importeRange.NumberFormat = "@"

$ws.Range('H2').Value = '4850.00'
$ws.Range('H3').Value = '4250.00'
$ws.Range('H4').Value = '39635.00'
$ws.Range('H5').Value = '260356.20'
$ws.Range('H6').Value = '1740.00'
$ws.Range('H7').Value = '999.65'
$ws.Range('H8').Value = '6957.50'
$ws.Range('H9').Value = '213750.00'
$ws.Range('H10').Value = '4099.95'
$ws.Range('H11').Value = '4200.16'
$ws.Range('H12').Value = '184894.21'
$ws.Range('H13').Value = '236499.91'
$ws.Range('H14').Value = '40680.00'
$ws.Range('H15').Value = '11238.20'
$ws.Range('H16').Value = '115371.84'
$ws.Range('H17').Value = '4850.00'
$ws.Range('H18').Value = '23400.00'
$ws.Range('H19').Value = '7479.46'
$ws.Range('H20').Value = '12880.00'
$ws.Range('H21').Value = '650.00'
$ws.Range('H22').Value = '4868.00'
$ws.Range('H23').Value = '8380.00'
$ws.Range('H24').Value = '576.00'
$ws.Range('H25').Value = '28909.99'
$ws.Range('H26').Value = '1511.20'
$ws.Range('H27').Value = '86.32'
$ws.Range('H28').Value = '50558.60'
$ws.Range('H29').Value = '43857.00'
$ws.Range('H30').Value = '40619.27'
$ws.Range('H31').Value = '720.00'
$ws.Range('H32').Value = '16351.00'
$ws.Range('H33').Value = '6926.73'
$ws.Range('H34').Value = '6500.00'
$ws.Range('H35').Value = '55.50'
$ws.Range('H36').Value = '50.31'
$ws.Range('H37').Value = '12490.00'
$ws.Range('H38').Value = '7705.00'
$ws.Range('H39').Value = '30466.04'
$ws.Range('H40').Value = '76.00'
$ws.Range('H41').Value = '5368.39'
$ws.Range('H42').Value = '38.00'
$ws.Range('H43').Value = '6752.00'
$ws.Range('H44').Value = '6827.15'
$ws.Range('H45').Value = '35977.22'
$ws.Range('H46').Value = '1553.30'
$ws.Range('H47').Value = '1956.66'
$ws.Range('H48').Value = '136348.00'
$ws.Range('H49').Value = '2623.82'
$ws.Range('H50').Value = '2967.42'
$ws.Range('H51').Value = '774.00'
$ws.Range('H52').Value = '1765.89'
$ws.Range('H53').Value = '6913.86'
$ws.Range('H54').Value = '9704.53'
$ws.Range('H55').Value = '5546.62'
$ws.Range('H56').Value = '42108.58'
$ws.Range('H57').Value = '95.40'
$ws.Range('H58').Value = '1653.24'
$ws.Range('H59').Value = '240.00'
$ws.Range('H60').Value = '35386.00'
$ws.Range('H61').Value = '2978.00'
$ws.Range('H62').Value = '20400.00'
$ws.Range('H63').Value = '270.94'
$ws.Range('H64').Value = '5926.65'
$ws.Range('H65').Value = '755.03'
$ws.Range('H66').Value = '3889.99'
$ws.Range('H67').Value = '11094.10'
$ws.Range('H68').Value = '2850.00'
$ws.Range('H69').Value = '10278.00'
$ws.Range('H70').Value = '6000.00'
$ws.Range('H71').Value = '2710.00'
$ws.Range('H72').Value = '39000.00'
$ws.Range('H73').Value = '1815.75'
$ws.Range('H74').Value = '4000.00'
$ws.Range('H75').Value = '9690.00'
$ws.Range('H76').Value = '1700.00'
$ws.Range('H77').Value = '300.00'
$ws.Range('H78').Value = '1600.00'
$ws.Range('H79').Value = '16886.94'
$ws.Range('H80').Value = '6594.00'
$ws.Range('H81').Value = '10320.00'
$ws.Range('H82').Value = '9900.00'
$ws.Range('H83').Value = '21150.00'
$ws.Range('H84').Value = '5850.00'
$ws.Range('H85').Value = '90.00'
$ws.Range('H86').Value = '7796.00'
$ws.Range('H87').Value = '141.00'
$ws.Range('H88').Value = '3905.00'
$ws.Range('H89').Value = '14081.00'
$ws.Range('H90').Value = '2525.99'
$ws.Range('H91').Value = '10000.00'
$ws.Range('H92').Value = '100.00'
$ws.Range('H93').Value = '1920.00'
$ws.Range('H94').Value = '14000.00'
$ws.Range('H95').Value = '45638.43'
$ws.Range('H96').Value = '11799.99'
$ws.Range('H97').Value = '8190.67'
$ws.Range('H98').Value = '345761.96'
$ws.Range('H99').Value = '350.00'
$ws.Range('H100').Value = '36474.25'
$ws.Range('H101').Value = '5710.00'
$ws.Range('H102').Value = '1800.00'
$ws.Range('H103').Value = '21.49'
$ws.Range('H104').Value = '3261.00'
$ws.Range('H105').Value = '0.63'
$ws.Range('H106').Value = '172.00'
$ws.Range('H107').Value = '200.00'
$ws.Range('H108').Value = '850.00'
$ws.Range('H109').Value = '23400.00'
$ws.Range('H110').Value = '20772.36'
$ws.Range('H111').Value = '38169.49'
$ws.Range('H112').Value = '4210.00'
$ws.Range('H113').Value = '19.87'
$ws.Range('H114').Value = '1800.00'
$ws.Range('H115').Value = '16337.30'
$ws.Range('H116').Value = '2819.00'
$ws.Range('H117').Value = '293.10'
$ws.Range('H118').Value = '1206.20'
$ws.Range('H119').Value = '78000.00'
$ws.Range('H120').Value = '110.01'
$ws.Range('H121').Value = '4130.00'
$ws.Range('H122').Value = '760.00'
$ws.Range('H123').Value = '309.00'
$ws.Range('H124').Value = '25535.00'
$ws.Range('H125').Value = '948.00'
$ws.Range('H126').Value = '360.00'
$ws.Range('H127').Value = '3580.00'
$ws.Range('H128').Value = '332.52'
$ws.Range('H129').Value = '4250.00'
$ws.Range('H130').Value = '2000.00'
$ws.Range('H131').Value = '750.00'
$ws.Range('H132').Value = '204750.00'
$ws.Range('H133').Value = '3700.00'
$ws.Range('H134').Value = '4200.00'
$ws.Range('H135').Value = '3533.20'
$ws.Range('H136').Value = '18650.00'
$ws.Range('H137').Value = '8500.00'
$ws.Range('H138').Value = '2310.00'
$ws.Range('H139').Value = '1695.00'
$ws.Range('H140').Value = '1284.10'
$ws.Range('H141').Value = '117976.00'
$ws.Range('H142').Value = '4000.00'
$ws.Range('H143').Value = '23500.00'
$ws.Range('H144').Value = '9600.00'
$ws.Range('H145').Value = '1000.00'
$ws.Range('H146').Value = '2285.00'
$ws.Range('H147').Value = '6000.00'
$ws.Range('H148').Value = '211467.00'
$ws.Range('H149').Value = '420877.00'
$ws.Range('H150').Value = '1692.06'
$ws.Range('H151').Value = '515.45'
$ws.Range('H152').Value = '330.50'
$ws.Range('H153').Value = '121.00'
$ws.Range('H154').Value = '400.00'
$ws.Range('H155').Value = '684360.00'
$ws.Range('H156').Value = '210270.52'
$ws.Range('H157').Value = '41938.54'
$ws.Range('H158').Value = '3300.00'
$ws.Range('H159').Value = '4300.00'
$ws.Range('H160').Value = '1500.00'
$ws.Range('H161').Value = '8137.25'
$ws.Range('H162').Value = '1200.00'
$ws.Range('H163').Value = '1500.00'
$ws.Range('H164').Value = '1800.00'
$ws.Range('H165').Value = '1000.00'
$ws.Range('H166').Value = '3000.00'
$ws.Range('H167').Value = '2300.00'
$ws.Range('H168').Value = '1300.00'
$ws.Range('H169').Value = '800.00'
$ws.Range('H170').Value = '500.00'
$ws.Range('H171').Value = '9280.00'
$ws.Range('H172').Value = '4500.00'
$ws.Range('H173').Value = '4200.00'
$ws.Range('H174').Value = '10080.00'
$ws.Range('H175').Value = '1300.00'
$ws.Range('H176').Value = '585.00'
$ws.Range('H177').Value = '50644.50'
$ws.Range('H178').Value = '500.00'
$ws.Range('H179').Value = '9650.00'
$ws.Range('H180').Value = '32.50'
$ws.Range('H181').Value = '2444.00'
$ws.Range('H182').Value = '980.00'
$ws.Range('H183').Value = '5300.00'
$ws.Range('H184').Value = '1455.95'
$ws.Range('H185').Value = '21720.00'
$ws.Range('H186').Value = '7015.00'
$ws.Range('H187').Value = '5474.24'
$ws.Range('H188').Value = '405.50'
$ws.Range('H189').Value = '2812.60'
$ws.Range('H190').Value = '744.50'
$ws.Range('H191').Value = '11920.00'
$ws.Range('H192').Value = '14399.96'
$ws.Range('H193').Value = '1713.00'
$ws.Range('H194').Value = '5330.61'
$ws.Range('H195').Value = '11897.62'
$ws.Range('H196').Value = '2360.00'
$ws.Range('H197').Value = '470.00'
$ws.Range('H198').Value = '1626.00'
$ws.Range('H199').Value = '2700.00'
$ws.Range('H200').Value = '630.00'
$ws.Range('H201').Value = '4295.20'
$ws.Range('H202').Value = '2415.20'
$ws.Range('H203').Value = '5295.93'
$ws.Range('H204').Value = '4800.00'
$ws.Range('H205').Value = '12200.00'
$ws.Range('H206').Value = '208120.00'
$ws.Range('H207').Value = '16800.00'
$ws.Range('H208').Value = '6292.00'
$ws.Range('H209').Value = '4500.00'
$ws.Range('H210').Value = '245520.00'
$ws.Range('H211').Value = '2744849.51'
$ws.Range('H212').Value = '129544.00'
$ws.Range('H213').Value = '92780.00'
$ws.Range('H214').Value = '120372.00'
$ws.Range('H215').Value = '128000.00'
$ws.Range('H216').Value = '6473830.54'
$ws.Range('H217').Value = '28154.50'
$ws.Range('H218').Value = '190000.00'
$ws.Range('H219').Value = '217993.70'
$ws.Range('H220').Value = '5488.40'
$ws.Range('H221').Value = '122990.00'
$ws.Range('H222').Value = '45000.00'
$ws.Range('H223').Value = '187400.00'
$ws.Range('H224').Value = '4000.00'
$ws.Range('H225').Value = '4435.60'
$ws.Range('H226').Value = '2100.00'
$ws.Range('H227').Value = '22400.00'
$ws.Range('H228').Value = '23120.00'
$ws.Range('H229').Value = '2784.95'
$ws.Range('H230').Value = '53500.00'
$ws.Range('H231').Value = '7200.00'

$importeRange.Style = "Normal"

# --- Fix comma-separated names (commas -> periods, stray periods dropped) ---
$ws.Range('E40').Value = 'ALBIZZATTI. PABLO MARTIN Y FULINI. SERGIO RUBEN'
$ws.Range('E183').Value = 'ALBIZZATTI. PABLO MARTIN Y FULINI. SERGIO RUBEN'
$ws.Range('E84').Value = 'MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO'
$ws.Range('E130').Value = 'RICCOTTI. MARIANA EDITH'
$ws.Range('E182').Value = 'SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH'
$ws.Range('E192').Value = 'SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH'
$ws.Range('F147').Value = 'MERCANZINI. GASTON ARIEL'

